$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info (rows 3-5) ---
$ws.Range("B3").Value = "12/16/2020"
$ws.Range("B4").Value = "Limette"
$ws.Range("B5").Value = 4

# --- Team member names + salaries (rows 8-11) ---
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("B8").Value = 100
$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("B9").Value = 100
$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("B10").Value = 100
$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("B11").Value = 100

# Row 12 previously held "Member 5" - now cleared
$ws.Range("A12").Value = $null

# --- Tasks header row shrinks ---
$ws.Rows.Item(18).RowHeight = 39

# --- Tasks completed / to complete next week (rows 19-22) ---
$ws.Range("A19").Value = "Finalized communication of front- and backend"
$ws.Range("B19").Value = "Clean up git repo"
$ws.Range("A20").Value = "Prepared final presentation"
$ws.Range("B20").Value = $null
$ws.Range("A21").Value = "Collected code build and run instructions"
$ws.Range("B21").Value = $null
$ws.Range("A22").Value = "Created video for hifi prototype"
$ws.Range("B22").Value = $null

# Apply the smaller 10pt font used for the new task rows while keeping
# the existing thin border around each cell.
$taskRange = $ws.Range("A19:B22")
$taskRange.Font.Size = 10
$taskRange.Font.Scheme = "minor"
$taskRange.Borders.LineStyle = 1

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moves to B4 ---
$ws.Range("B4").Select() | Out-Null
